# The commit swaps the content of ppt/theme/theme1.xml ("Office Theme")
# and ppt/theme/theme2.xml ("Integral") -- i.e. the deck's active theme
# (wired to the one SlideMaster / the whole presentation) changes from the
# "Integral" color palette to the "Office Theme" color palette.
#
# theme1.xml and theme2.xml share an identical <a:fontScheme> and
# <a:fmtScheme>; the only real difference between the two theme parts is
# the <a:clrScheme> (12 RGB slots) plus the cosmetic name="" attributes.
# The PowerPoint object model exposes exactly that color scheme via
# Master.ColorScheme.Colors(1..12).RGB, so drive the swap through it.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.ColorScheme

function ToRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $b * 65536 + $g * 256 + $r
}

# Office Theme palette (was theme1.xml), in clrScheme slot order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $cs.Colors($i).RGB = ToRGB($officeThemeColors[$i - 1])
}
